$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells ---
# Same column layout as before (Colour->B, Clarity->C, Cut->D,
# Certification->E, Polish->F, Symmetry->G), just renamed to the
# "categorical_*" labels; A1 (Carat) and H1 (price) keep their text.
$ws.Range("A1").Value2 = "Carat"
$ws.Range("B1").Value2 = "categorical_color"
$ws.Range("C1").Value2 = "categorical_clarity"
$ws.Range("D1").Value2 = "categorical_cut"
$ws.Range("E1").Value2 = "categorical_certification"
$ws.Range("F1").Value2 = "categorical_polish"
$ws.Range("G1").Value2 = "categorical_symmetry"
$ws.Range("H1").Value2 = "price"

# --- Re-style the "categorical_*" header cells ---
# D1 already carried the bold-Arial-10 style (fontId 2 / xf 2); restyle that
# font in place (Courier New 11, black, non-bold) and add left/center
# alignment, then fan the resulting format out to the other new columns.
$d1 = $ws.Range("D1")
$d1.Font.Name = "Courier New"
$d1.Font.Size = 11
$d1.Font.Bold = $false
$d1.Font.Color = 0
$d1.Font.Family = 3
$d1.HorizontalAlignment = -4131
$d1.VerticalAlignment = -4108

[void]$d1.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Page setup: force portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Move the active selection to C1 ---
[void]$ws.Range("C1").Select()
